# The deck's slide master (ppt/theme/theme1.xml, theme name "Integral") and
# notes master (ppt/theme/theme2.xml, theme name "Office Theme") swap their
# full theme content (color scheme + theme name). The fontScheme/fmtScheme
# are identical between the two themes already, so the effective content
# swap is carried entirely by the 12 color-scheme entries (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink).
#
# Apply the target "Office Theme" palette (currently living in theme2.xml)
# onto the presentation's editable theme color scheme, via the PowerPoint
# object model (Master.ColorScheme), which is the channel the host exposes
# for rewriting ppt/theme/theme1.xml's <a:clrScheme>.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

# ppColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink. RGBColor.RGB is the usual VBA BGR-packed value,
# i.e. for hex RRGGBB the stored long is (BB*65536 + GG*256 + RR).
function Set-SchemeColor($index, $r, $g, $b) {
    $cs.Colors($index).RGB = ($b * 65536) + ($g * 256) + $r
}

Set-SchemeColor 1  0x00 0x00 0x00   # dk1      000000
Set-SchemeColor 2  0xFF 0xFF 0xFF   # lt1      FFFFFF
Set-SchemeColor 3  0x44 0x54 0x6A   # dk2      44546A
Set-SchemeColor 4  0xE7 0xE6 0xE6   # lt2      E7E6E6
Set-SchemeColor 5  0x5B 0x9B 0xD5   # accent1  5B9BD5
Set-SchemeColor 6  0xED 0x7D 0x31   # accent2  ED7D31
Set-SchemeColor 7  0xA5 0xA5 0xA5   # accent3  A5A5A5
Set-SchemeColor 8  0xFF 0xC0 0x00   # accent4  FFC000
Set-SchemeColor 9  0x44 0x72 0xC4   # accent5  4472C4
Set-SchemeColor 10 0x70 0xAD 0x47   # accent6  70AD47
Set-SchemeColor 11 0x05 0x63 0xC1   # hlink    0563C1
Set-SchemeColor 12 0x95 0x4F 0x72   # folHlink 954F72
